$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.497.63"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").Value = "1.830.88"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.97"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4287"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07273"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8690"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.65"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").Value = "1.827.90"
$ws.Range("E12").Value = "  +1.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.409"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.534"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06938"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.37"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008892"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.41"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").Value = "27.795.89"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.128"
$ws.Range("E22").Value = "  +3.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("E23").Value = "  +4.47%  "

$ws.Range("D24").Value = "2.103.22"
$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.38"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.146"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.14"
$ws.Range("E29").Value = "  -5.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.834"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08881"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7576"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.984"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.538"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.095"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05317"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.805"
$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5076"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1661"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.617"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.394"
$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.49"
$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.05"
$ws.Range("E46").Value = "  +1.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06503"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4683"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.612"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.746"
$ws.Range("E51").Value = "  +2.78%  "
